# LOQ4093.docx edit: the section blocks got reshuffled down by one slot
# (with the last block wrapping back to the "Docente(s)" slot). No
# paragraphs are added or removed and no paragraph styles change - only
# the text content moves between the already-existing paragraphs/runs.
#
# [char]11 is Word's "manual line break" character (rendered as <w:br/>
# in the XML); we splice it in wherever the source run was followed by
# a <w:br/>. Range.Text includes the trailing paragraph mark ([char]13),
# which must be trimmed off before reusing the captured string as a
# replacement, otherwise assigning it back in splits the paragraph.

$d = $word.ActiveDocument
$cr = [char]13

function Get-ParaText($index) {
    return $d.Paragraphs($index).Range.Text.TrimEnd($cr)
}

# --- Snapshot all the "old" text we are about to move, before any of it
# --- gets overwritten (the moves form a single rotation cycle). ---

$oldObjetivosBody    = Get-ParaText 6     # "Visao integrada..."
$oldDocenteNome      = Get-ParaText 8     # "1285870 - Marcos Villela Barcza"
$oldProgResumido     = Get-ParaText 10    # "1.Petroleo: ... 10.Oleos ..."
$oldProgramaBody     = Get-ParaText 12    # "Petroleo: historico; ... Oleos..."
$oldBibliografiaBody = Get-ParaText 16    # "a)Speight ... f)Revista Petro & Quimica."

$avalPara = $d.Paragraphs(14)

$oldMetodoText = $null
$avalRange = $avalPara.Range.Duplicate()
if ($avalRange.Find.Execute("Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos e seminários")) {
    $oldMetodoText = $avalRange.Text
}

$oldCriterioText = $null
$avalRange = $avalPara.Range.Duplicate()
if ($avalRange.Find.Execute("Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula.")) {
    $oldCriterioText = $avalRange.Text
}

$oldNormaText = $null
$avalRange = $avalPara.Range.Duplicate()
if ($avalRange.Find.Execute("Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação.")) {
    $oldNormaText = $avalRange.Text
}

# --- Now write the rotated content back in, slot by slot. ---

# 1) "Objetivos" body paragraph <- old "Programa resumido" body (the long
#    numbered list with manual line breaks).
$d.Paragraphs(6).Range.Text = $oldProgResumido

# 2) "Docente(s) Responsavel(eis)" list item <- old "Objetivos" body text.
$d.Paragraphs(8).Range.Text = $oldObjetivosBody

# 3) "Programa resumido" body paragraph <- old "Programa" body (the
#    semicolon-joined single-paragraph version).
$d.Paragraphs(10).Range.Text = $oldProgramaBody

# 4) "Programa" body paragraph <- old "Metodo:" value text.
$d.Paragraphs(12).Range.Text = $oldMetodoText

# 5) Avaliacao paragraph: shift each value one label down, and the last
#    one ("Norma de recuperacao:") picks up the old Bibliografia list.
$avalRange = $avalPara.Range.Duplicate()
$avalRange.Find.Execute($oldMetodoText, $true, $false, $false, $false, $false, $true, 1, $false, $oldCriterioText, 2) | Out-Null

$avalRange = $avalPara.Range.Duplicate()
$avalRange.Find.Execute($oldCriterioText, $true, $false, $false, $false, $false, $true, 1, $false, $oldNormaText, 2) | Out-Null

$avalRange = $avalPara.Range.Duplicate()
$avalRange.Find.Execute($oldNormaText, $true, $false, $false, $false, $false, $true, 1, $false, $oldBibliografiaBody, 2) | Out-Null

# 6) Bibliografia body paragraph <- old Docente(s) name ("1285870 - ...").
$d.Paragraphs(16).Range.Text = $oldDocenteNome
